$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "[8, 35, 45, 8, 5]"
$ws.Range("B7").Value = 0.9732475505600137
$ws.Range("B8").Value = 0.01324315225022732
$ws.Range("B9").Value = 29
$ws.Range("B10").Value = 80
$ws.Range("B12").Value = "[[28, 39], [63, 68]]"
$ws.Range("B14").Value = "[[79, 76, 4], [83, 47, 88], [35, 25, 96]]"
$ws.Range("B15").Value = "[0.723570428093315, 0.8243709715886228, 0.9360453142502849, 0.8248084839896768, 0.7532945172554409, 0.6654888335215752]"
